$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right before the existing "MegaSync" row (row 79),
# shifting every row below it down by one.
$ws.Rows("79:79").Insert() | Out-Null

# Copy the formatting of the row that just got pushed down (now row 80)
# into the freshly inserted, still-blank row 79 so the new row keeps the
# same borders/fill/font as the rest of the table.
$ws.Range("B80:G80").Copy() | Out-Null
$ws.Range("B79:G79").PasteSpecial(-4122) | Out-Null

# Populate the new row with the new alias entry.
$ws.Range("C79").Value = "Librewolf (Mails)"
$ws.Range("D79").Value = "mails"
$ws.Range("E79").Value = "[Application]"

# The row-number column (B3:B159) is a dynamic SEQUENCE() array formula;
# re-enter it so its spill range grows to cover the newly added row.
$ws.Range("B3:B160").FormulaArray = "=SEQUENCE(COUNTA(C:C)-1)"

# Restore focus roughly where the author left it.
$ws.Range("D87").Select()
